# Adds a new "tenant_id" column to the org.xlsx export template, between
# the existing "update_time" columns' predecessor (update_usr_id) and the
# rest, i.e. right after "create_time" / before "update_usr_id":
#
#   ... | create_usr_id | create_time | tenant_id | update_usr_id | update_time
#
# Row 1 holds the column comment/header template text, row 2 holds the
# per-row model template text (the "forRow" data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "update_usr_id_lbl" / "update_time_lbl" columns
# (H, I) one column to the right (-> I, J), opening up column H for the
# new tenant_id column.
$ws.Columns("H").Insert()

# Row 1 (comments): header + data-validation-list template for the new
# tenant_id column.
$ws.Range("H1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

# Row 2 (model/data row template): the tenant_id label cell.
$ws.Range("H2").Value = '<%=model.tenant_id_lbl%>'
